# Update the arithmetic-practice table: every one of the 100 problem
# cells (20 rows x 5 columns) is replaced with a new "a+b=" / "a-b="
# expression per the target revision. Cells are addressed positionally
# (row, column) via Table.Cell(r, c) rather than by searching for the
# old text, since several expressions (e.g. "3+63=") repeat verbatim
# in the original table and a text-based Find/Replace would be
# ambiguous.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "30+18="
$t.Cell(1,2).Range.Text = "51+14="
$t.Cell(1,3).Range.Text = "95-69="
$t.Cell(1,4).Range.Text = "62-11="
$t.Cell(1,5).Range.Text = "44+11="
$t.Cell(2,1).Range.Text = "86-52="
$t.Cell(2,2).Range.Text = "63-45="
$t.Cell(2,3).Range.Text = "86-45="
$t.Cell(2,4).Range.Text = "40+31="
$t.Cell(2,5).Range.Text = "95-89="
$t.Cell(3,1).Range.Text = "25-14="
$t.Cell(3,2).Range.Text = "52+25="
$t.Cell(3,3).Range.Text = "54-26="
$t.Cell(3,4).Range.Text = "99-13="
$t.Cell(3,5).Range.Text = "27-6="
$t.Cell(4,1).Range.Text = "58-52="
$t.Cell(4,2).Range.Text = "43+51="
$t.Cell(4,3).Range.Text = "60-9="
$t.Cell(4,4).Range.Text = "49+46="
$t.Cell(4,5).Range.Text = "49+11="
$t.Cell(5,1).Range.Text = "58+40="
$t.Cell(5,2).Range.Text = "32+25="
$t.Cell(5,3).Range.Text = "45-25="
$t.Cell(5,4).Range.Text = "86-10="
$t.Cell(5,5).Range.Text = "94-22="
$t.Cell(6,1).Range.Text = "94-16="
$t.Cell(6,2).Range.Text = "26+25="
$t.Cell(6,3).Range.Text = "32+53="
$t.Cell(6,4).Range.Text = "60-0="
$t.Cell(6,5).Range.Text = "66-31="
$t.Cell(7,1).Range.Text = "54-37="
$t.Cell(7,2).Range.Text = "79-48="
$t.Cell(7,3).Range.Text = "19+16="
$t.Cell(7,4).Range.Text = "83-6="
$t.Cell(7,5).Range.Text = "90-42="
$t.Cell(8,1).Range.Text = "20+10="
$t.Cell(8,2).Range.Text = "82-63="
$t.Cell(8,3).Range.Text = "47+20="
$t.Cell(8,4).Range.Text = "79+4="
$t.Cell(8,5).Range.Text = "17+51="
$t.Cell(9,1).Range.Text = "63-35="
$t.Cell(9,2).Range.Text = "60-10="
$t.Cell(9,3).Range.Text = "61+26="
$t.Cell(9,4).Range.Text = "39+18="
$t.Cell(9,5).Range.Text = "43-18="
$t.Cell(10,1).Range.Text = "50-49="
$t.Cell(10,2).Range.Text = "81-77="
$t.Cell(10,3).Range.Text = "35-13="
$t.Cell(10,4).Range.Text = "11+50="
$t.Cell(10,5).Range.Text = "81-69="
$t.Cell(11,1).Range.Text = "82-49="
$t.Cell(11,2).Range.Text = "93-34="
$t.Cell(11,3).Range.Text = "33-24="
$t.Cell(11,4).Range.Text = "37-27="
$t.Cell(11,5).Range.Text = "31+9="
$t.Cell(12,1).Range.Text = "25+4="
$t.Cell(12,2).Range.Text = "63-19="
$t.Cell(12,3).Range.Text = "14+37="
$t.Cell(12,4).Range.Text = "87-76="
$t.Cell(12,5).Range.Text = "60+18="
$t.Cell(13,1).Range.Text = "53+28="
$t.Cell(13,2).Range.Text = "73+7="
$t.Cell(13,3).Range.Text = "39+27="
$t.Cell(13,4).Range.Text = "15+71="
$t.Cell(13,5).Range.Text = "82-15="
$t.Cell(14,1).Range.Text = "59-47="
$t.Cell(14,2).Range.Text = "75-59="
$t.Cell(14,3).Range.Text = "63+34="
$t.Cell(14,4).Range.Text = "40-8="
$t.Cell(14,5).Range.Text = "91-22="
$t.Cell(15,1).Range.Text = "21+78="
$t.Cell(15,2).Range.Text = "97-41="
$t.Cell(15,3).Range.Text = "75-10="
$t.Cell(15,4).Range.Text = "46+52="
$t.Cell(15,5).Range.Text = "20+20="
$t.Cell(16,1).Range.Text = "20+2="
$t.Cell(16,2).Range.Text = "80+9="
$t.Cell(16,3).Range.Text = "47+15="
$t.Cell(16,4).Range.Text = "89-74="
$t.Cell(16,5).Range.Text = "90-83="
$t.Cell(17,1).Range.Text = "72-70="
$t.Cell(17,2).Range.Text = "68-35="
$t.Cell(17,3).Range.Text = "32-8="
$t.Cell(17,4).Range.Text = "94-38="
$t.Cell(17,5).Range.Text = "59+19="
$t.Cell(18,1).Range.Text = "88-65="
$t.Cell(18,2).Range.Text = "65-62="
$t.Cell(18,3).Range.Text = "50+32="
$t.Cell(18,4).Range.Text = "27+8="
$t.Cell(18,5).Range.Text = "65+12="
$t.Cell(19,1).Range.Text = "59-21="
$t.Cell(19,2).Range.Text = "62+36="
$t.Cell(19,3).Range.Text = "82-8="
$t.Cell(19,4).Range.Text = "78+8="
$t.Cell(19,5).Range.Text = "43+53="
$t.Cell(20,1).Range.Text = "64+17="
$t.Cell(20,2).Range.Text = "56-23="
$t.Cell(20,3).Range.Text = "40+0="
$t.Cell(20,4).Range.Text = "97-28="
$t.Cell(20,5).Range.Text = "29+59="
